$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.174.62"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "2.269.25"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("B5").Value = "BinanceUSD"
$ws.Range("C5").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D5").Value = "'106.18"
$ws.Range("E5").Value = "  +10,506.06%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'305.92"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'93.58"
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.532"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").Value = "'33.05"
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").Value = "'0.112"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "2.621.55"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "'14.38"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").Value = "2.268.30"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "'0.787"
$ws.Range("E18").Value = "  +4.00%  "
$ws.Range("D19").Value = "42.028.28"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "'12.72"
$ws.Range("E20").Value = "  +5.24%  "
$ws.Range("D21").Value = "0.0₃0919"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("D23").Value = "'68.22"
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("D24").Value = "'244.19"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("E25").Value = "  +2.16%  "
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D31").Value = "'35.18"
$ws.Range("E31").Value = "  +3.95%  "
$ws.Range("D32").Value = "'159.89"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'0.0745"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("E37").Value = "  +3.81%  "
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("E42").Value = "  +3.66%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'19.67"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.011.51"
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("E45").Value = "  +10.46%  "
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("D47").Value = "'10.24"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("D49").Value = "'53.86"
$ws.Range("E49").Value = "  +4.11%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "'72.82"
$ws.Range("E50").Value = "  +2.79%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.52"
$ws.Range("E51").Value = "  +0.59%  "
